$wb = $excel.ActiveWorkbook

# Sheet "식당판매" (restaurant sales): 수저(spoon) count 3 -> 2
$ws1 = $wb.Worksheets.Item("식당판매")
$ws1.Range("C6").Value = 2

# Sheet "상복" (mourning clothes): 식탁보(tablecloth) count 3 -> 2
$ws4 = $wb.Worksheets.Item("상복")
$ws4.Range("C7").Value = 2

# Sheet "기타" (other): updated counts from center
$ws5 = $wb.Worksheets.Item("기타")
$ws5.Range("C8").Value = 39
$ws5.Range("C9").Value = 11
$ws5.Range("C10").Value = 77
$ws5.Range("C11").Value = 20
